# WRI Updates to HK model from 11/15
$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsRM = $wb.Worksheets.Item("RM")

# --- About sheet: rewrite the Notes text block ---
$wsAbout.Range("A10").Value = "The reserve margin in the U.S. dataset doesn't vary by year, but the RM Reserve Margin variable is a time series to support"
$wsAbout.Range("A11").Value = "countries that project changes in future reserve margin by year."
$wsAbout.Range("A12").Value = $null
$wsAbout.Range("A13").Value = "use China variables"
$wsAbout.Range("A13").Font.Name = "宋体"
$wsAbout.Range("A13").Font.Color = 255

# --- RM sheet: update reserve margin values and formatting ---
$wsRM.Range("A1").Value = ""
$rmRange = $wsRM.Range("A1:AK2")
$rmRange.WrapText = $true
$rmRange.VerticalAlignment = -4108

$wsRM.Range("B2:AK2").Value = 0.15

# --- Selection / active sheet state ---
$wsAbout.Range("A13").Select()
$wsRM.Range("A1:AK2").Select()
$wsRM.Activate()
